$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price strings stay as text (preserve trailing zeros / exact formatting)
$textCells = @("D5", "D6", "D9", "D10", "D12", "D19", "D20", "D21", "D22", "D23", "D24", "D29", "D31", "D32", "D35", "D37", "D38", "D39", "D40", "D42", "D43", "D44", "D45", "D47", "D50", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply cell value updates from the diff
$ws.Range("D2").Value = "62.787.88"
$ws.Range("E2").Value = "  +1.86%  "
$ws.Range("D3").Value = "3.459.19"
$ws.Range("E3").Value = "  +1.81%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "577.65"
$ws.Range("E5").Value = "  +0.27%  "
$ws.Range("D6").Value = "146.63"
$ws.Range("E6").Value = "  +3.14%  "
$ws.Range("D7").Value = "3.459.90"
$ws.Range("E7").Value = "  +1.80%  "
$ws.Range("E8").Value = "  -0.07%  "
$ws.Range("D9").Value = "0.480"
$ws.Range("E9").Value = "  +1.58%  "
$ws.Range("D10").Value = "7.64"
$ws.Range("E10").Value = "  -0.06%  "
$ws.Range("E11").Value = "  +0.84%  "
$ws.Range("D12").Value = "0.400"
$ws.Range("E12").Value = "  +3.59%  "
$ws.Range("D13").Value = "4.049.75"
$ws.Range("E13").Value = "  +1.82%  "
$ws.Range("E14").Value = "  +5.14%  "
$ws.Range("E15").Value = "  +2.85%  "
$ws.Range("D16").Value = "3.491.64"
$ws.Range("E16").Value = "  +2.91%  "
$ws.Range("E17").Value = "  -0.13%  "
$ws.Range("D18").Value = "62.753.36"
$ws.Range("E18").Value = "  +1.66%  "
$ws.Range("D19").Value = "6.29"
$ws.Range("E19").Value = "  +2.68%  "
$ws.Range("D20").Value = "14.25"
$ws.Range("E20").Value = "  +4.46%  "
$ws.Range("D21").Value = "9.16"
$ws.Range("E21").Value = "  +1.52%  "
$ws.Range("D22").Value = "387.56"
$ws.Range("E22").Value = "  -0.29%  "
$ws.Range("B23").Value = "Litecoin"
$ws.Range("C23").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D23").Value = "74.41"
$ws.Range("E23").Value = "  -0.31%  "
$ws.Range("B24").Value = "Polygon"
$ws.Range("C24").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D24").Value = "0.554"
$ws.Range("E24").Value = "  +0.99%  "
$ws.Range("E25").Value = "  -0.03%  "
$ws.Range("D26").Value = "3.606.03"
$ws.Range("E26").Value = "  +2.02%  "
$ws.Range("E27").Value = "  +0.45%  "
$ws.Range("E28").Value = "  -9.16%  "
$ws.Range("D29").Value = "7.49"
$ws.Range("E29").Value = "  +1.32%  "
$ws.Range("D31").Value = "8.12"
$ws.Range("E31").Value = "  +1.49%  "
$ws.Range("D32").Value = "2.12"
$ws.Range("E32").Value = "  -1.19%  "
$ws.Range("E33").Value = "  +0.06%  "
$ws.Range("E34").Value = "  -4.31%  "
$ws.Range("D35").Value = "23.56"
$ws.Range("E35").Value = "  +0.89%  "
$ws.Range("E36").Value = "  +3.66%  "
$ws.Range("D37").Value = "7.01"
$ws.Range("E37").Value = "  +1.36%  "
$ws.Range("D38").Value = "31.38"
$ws.Range("E38").Value = "  +18.78%  "
$ws.Range("B39").Value = "Monero"
$ws.Range("C39").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D39").Value = "169.90"
$ws.Range("E39").Value = "  +0.77%  "
$ws.Range("B40").Value = "ImmutableX"
$ws.Range("C40").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D40").Value = "1.56"
$ws.Range("E40").Value = "  +5.60%  "
$ws.Range("D41").Value = "3.495.86"
$ws.Range("E41").Value = "  +1.82%  "
$ws.Range("D42").Value = "0.0750"
$ws.Range("E42").Value = "  -2.12%  "
$ws.Range("D43").Value = "0.797"
$ws.Range("E43").Value = "  +2.24%  "
$ws.Range("D44").Value = "42.33"
$ws.Range("E44").Value = "  -0.44%  "
$ws.Range("D45").Value = "4.45"
$ws.Range("E45").Value = "  +0.42%  "
$ws.Range("E46").Value = "  +2.85%  "
$ws.Range("D47").Value = "1.19"
$ws.Range("E47").Value = "  +3.07%  "
$ws.Range("D48").Value = "2.594.45"
$ws.Range("E48").Value = "  +5.67%  "
$ws.Range("E49").Value = "  +10.33%  "
$ws.Range("D50").Value = "22.75"
$ws.Range("E50").Value = "  -0.16%  "
$ws.Range("D51").Value = "6.70"
$ws.Range("E51").Value = "  +0.60%  "
